$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 8: "Date" property value was refreshed to a newer export timestamp.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Row 20: "Case Sensitive" property now has its value populated with "true".
# A plain ".Value = 'true'" gets auto-coerced to a Boolean by the engine,
# so we build it as a text formula result first and then paste it back as
# a static value, which keeps the cell a genuine text ("s") cell.
$caseSensitiveCell = $ws.Range("B20")
$caseSensitiveCell.Formula = "=""true"""
$caseSensitiveCell.Copy()
$caseSensitiveCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
